$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 13 by copying row 11 (keeps identical cell styles:
# s="9" for A/B, s="1" for C/D) then overwrite its values below.
$ws.Rows("11").Copy()
$ws.Rows("13").Insert()

# Update the existing assert rows (11 and 12) with the new sample data.
$ws.Range("B11").Value = "667299000"

$ws.Range("B12").Value = "667299000"
$ws.Range("C12").Value = "3016875893"
$ws.Range("D12").Value = "732111198172290"

# Fill in the newly-inserted row 13.
$ws.Range("A13").Value = "10960370"
$ws.Range("B13").Value = "667299000"
$ws.Range("C13").Value = "3016875982"
$ws.Range("D13").Value = "732111198172291"

# Reflect the new active cell / selection (matches the saved view state).
$ws.Range("D14").Select()
